# Auto-generated edit script applying cell value updates per the target diff.
# All target cells are plain text (coin names, URLs, price/volume strings),
# so we force NumberFormat to Text before assigning values to avoid Excel's
# automatic number/percentage coercion, then restore the original style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '310.54'),
    @('E2', '-0.13%'),
    @('D3', '37.56'),
    @('E3', '-1.73%'),
    @('D4', '5.082'),
    @('E4', '-0.89%'),
    @('D5', '0.07759'),
    @('E5', '-3.94%'),
    @('D6', '4.359'),
    @('E6', '-1.87%'),
    @('D7', '8.222'),
    @('E7', '-0.89%'),
    @('D8', '1.891'),
    @('E8', '-2.76%'),
    @('B9', 'BTSEToken'),
    @('C9', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'),
    @('D9', '2.887'),
    @('E9', '-11.53%'),
    @('B10', 'MXToken'),
    @('C10', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D10', '0.9223'),
    @('E10', '-1.71%'),
    @('B11', 'LiechtensteinCryptoassetsExchange'),
    @('C11', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D11', '0.1212'),
    @('E11', '-8.06%'),
    @('B12', 'WazirX'),
    @('C12', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D12', '0.1916'),
    @('E12', '-1.05%'),
    @('B13', 'MandalaExchangeToken'),
    @('C13', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D13', '0.09268'),
    @('E13', '1.96%'),
    @('B14', 'BitrueCoin'),
    @('C14', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D14', '0.03432'),
    @('E14', '-1.49%'),
    @('B15', 'BitMartToken'),
    @('C15', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D15', '0.09677'),
    @('E15', '0.01%'),
    @('B16', 'BitForexToken'),
    @('C16', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D16', '0.001380'),
    @('E16', '-2.21%'),
    @('B17', 'TigerCash'),
    @('C17', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D17', '0.005858'),
    @('E17', '-4.62%'),
    @('B18', 'LEO'),
    @('C18', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D18', '3.551'),
    @('E18', '-0.62%'),
    @('E19', '-1.89%'),
    @('D20', '5.304'),
    @('E20', '5.57%'),
    @('D21', '0.1296'),
    @('E21', '1.10%'),
    @('D22', '0.2595'),
    @('E22', '3.90%'),
    @('D23', '0.02106'),
    @('E23', '5,595.40%'),
    @('D24', '0.04362'),
    @('E24', '0.10%'),
    @('D25', '0.001215'),
    @('E25', '-2.11%'),
    @('D26', '0.004251'),
    @('E26', '-10.00%'),
    @('D27', '0.0001303'),
    @('E27', '-65.59%'),
    @('D39', '0.02087'),
    @('E39', '-5.66%'),
    @('D40', '0.04975'),
    @('E40', '-5.10%'),
    @('D41', '0.007685'),
    @('E41', '1.25%'),
    @('D42', '0.009850'),
    @('E42', '-4.42%'),
    @('D43', '0.1345'),
    @('E43', '-3.05%'),
    @('D44', '0.001995'),
    @('E44', '-2.04%'),
    @('D45', '0.008838'),
    @('E45', '-2.85%'),
    @('D46', '0.00006674'),
    @('E46', '1.05%'),
    @('D47', '0.00000000752'),
    @('E47', '0.09%'),
    @('D48', '0.002939'),
    @('E48', '-2.42%'),
    @('D49', '0.001203'),
    @('E49', '-28.89%'),
    @('D50', '0.00002104'),
    @('E50', '0.09%'),
    @('D51', '0.0002004'),
    @('E51', '0.09%')
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}
